# Update the "tied_teams" column (O) so the team ordering matches the
# corrected ordering used right before the final matchday of the group stage.
# Only the team order within each list changes - the underlying sets of teams
# stay the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 3-4: ['Uruguay', 'Hungary'] -> ['Hungary', 'Uruguay']
for ($r = 3; $r -le 4; $r++) {
    $ws.Range("O$r").Value = "['Hungary', 'Uruguay']"
}

# Rows 40-52: ['Costa Rica', 'Ireland'] -> ['Ireland', 'Costa Rica']
for ($r = 40; $r -le 52; $r++) {
    $ws.Range("O$r").Value = "['Ireland', 'Costa Rica']"
}

# Rows 53-60: ['Colombia', 'Costa Rica', 'Ireland', 'Argentina'] -> ['Ireland', 'Argentina', 'Colombia', 'Costa Rica']
for ($r = 53; $r -le 60; $r++) {
    $ws.Range("O$r").Value = "['Ireland', 'Argentina', 'Colombia', 'Costa Rica']"
}

# Rows 61-62: ['Colombia', 'Argentina'] -> ['Argentina', 'Colombia']
for ($r = 61; $r -le 62; $r++) {
    $ws.Range("O$r").Value = "['Argentina', 'Colombia']"
}

# Rows 63-73: ['Colombia', 'Scotland', 'Austria', 'Argentina'] -> ['Austria', 'Scotland', 'Argentina', 'Colombia']
for ($r = 63; $r -le 73; $r++) {
    $ws.Range("O$r").Value = "['Austria', 'Scotland', 'Argentina', 'Colombia']"
}
